# Applies the "Added react native to skill set" edit to the SKILLS section:
#   - Programming Languages: "JavaScript, PHP"      -> "JavaScript, PHP 5.4+"
#   - Web Technologies     : "React, Node"           -> "Node, React, React Native"
#   - Databases             : "MySQL, SQL Server+, " -> "MySQL, " (and SQL Server+ moved to the end,
#                              after Firebase): "..., Redis, Firebase" -> "..., Firebase, , SQL Server+"
#
# The simulated Word OM coalesces any two (or more) contiguously adjacent runs
# that end up with byte-identical <w:rPr> whenever a Range.Text edit touches
# their boundary. To reproduce the run-split structure the diff expects, we
# therefore (1) rewrite the whole affected span's text in one assignment, and
# then (2) re-introduce each needed run boundary by toggling Bold on/off
# (back to its original value) on the trailing sub-range of that boundary --
# a formatting no-op that nonetheless forces the engine to split the run,
# re-creating separate <w:r> elements that each keep the correct, full
# <w:rPr>.

function Split-RunAt {
    param($doc, [int]$pos, [int]$endPos)
    $r = $doc.Range($pos, $endPos)
    $r.Font.Bold = $true
    $r2 = $doc.Range($pos, $endPos)
    $r2.Font.Bold = $false
}

function Replace-Segment {
    param($doc, $para, [string]$oldSeg, [string]$newSeg, [int[]]$relBounds)

    $full = $para.Range.Text
    $idx = $full.IndexOf($oldSeg)
    if ($idx -lt 0) {
        throw "Segment not found: $oldSeg"
    }
    $start = $para.Range.Start + $idx
    $end = $start + $oldSeg.Length

    $sub = $doc.Range($start, $end)
    $sub.Text = $newSeg

    $newEnd = $start + $newSeg.Length
    for ($i = $relBounds.Length - 1; $i -ge 0; $i--) {
        $pos = $start + $relBounds[$i]
        Split-RunAt $doc $pos $newEnd
    }
}

$d = $word.ActiveDocument

# --- Programming Languages: "JavaScript, PHP" -> "JavaScript" | ", PHP 5.4+" ---
# (boundary 0 re-separates the new text from the preceding ", " run, which
# would otherwise get swallowed into the edited span's merge)
$pProgLang = $d.Paragraphs.Item(36)
Replace-Segment $d $pProgLang "JavaScript, PHP" "JavaScript, PHP 5.4+" @(0, 10)

# --- Web Technologies: "React" | ", Node" -> "Node" | ", React, React Native" ---
$pWebTech = $d.Paragraphs.Item(37)
Replace-Segment $d $pWebTech "React, Node" "Node, React, React Native" @(0, 4)

# --- Databases: "MySQL, SQL Server+, " ... ", Redis" | ", Firebase"
#     -> "MySQL" | ", " ... "," | " " | "Firebase" | ", " | ", SQL Server+" ---
$pDatabases = $d.Paragraphs.Item(38)
Replace-Segment $d $pDatabases `
    "MySQL, SQL Server+, MongoDB, SQLite, Redis, Firebase" `
    "MySQL, MongoDB, SQLite, Firebase, , SQL Server+" `
    @(5, 7, 14, 22, 23, 24, 32, 34)

Write-Output "Programming Languages: $($pProgLang.Range.Text)"
Write-Output "Web Technologies     : $($pWebTech.Range.Text)"
Write-Output "Databases            : $($pDatabases.Range.Text)"
